$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.177.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.740.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.738.98'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.41'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000246'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.367.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.735.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.192.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.113'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +17.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.725'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000150'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.54%  '
$ws.Range("E31").Value = '  +6.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.887.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.674.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.87'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.323'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '434.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.26%  '
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("E45").Value = '  +2.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.776.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0354'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.76%  '
